# Completed Baseline judgmental forecasting analysis
#
# Swap the "naiveAR2"/"ifoCast" column labels (D1/E1) and their paired
# error-column labels (F1/G1 -> previously naiveAR2/ifoCast error labels,
# now G1/H1), and flip the sign of every computed error value in
# columns F:H for rows 2-47.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the header labels for the naiveAR2 / ifoCast columns ---
$ws.Range("D1").Value = "ifoCast"
$ws.Range("E1").Value = "naiveAR2"
$ws.Range("G1").Value = "error_realized_minus_ifoCast"
$ws.Range("H1").Value = "error_realized_minus_naiveAR2"

# --- Flip the sign of the error columns (F, G, H) for data rows 2-47 ---
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 47) { $lastRow = 47 }

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($col in @("F", "G", "H")) {
        $cell = $ws.Range("$col$r")
        $v = $cell.Value2
        if ($v -ne $null) {
            $cell.Value = -1 * $v
        }
    }
}
